$d = $word.ActiveDocument

$d.Content.Find.Execute("6. Korisnik klikće dugme “Ukloni”", $true, $false, $false, $false, $false,
                         $true, 1, $false, "6. Korisnik klikće dugme “Ukloni komentar”", 2)
